$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GCN_LSTM")

# 1. Rename worksheet to reflect PM10-specific results
$ws.Name = "GCN_LSTM_PM10"

# 2. Populate the new shared strings in the exact order that reproduces the
#    target shared-string table ordering (new text first appears this way
#    while already-referenced strings 0-24 are left completely untouched).
$ws.Range("A8").Value = "Waste Facilities"
$ws.Range("D2").Value = "Uniform random values between [0, 1]."
$ws.Range("D3").Value = "Population difference across disctricts."
$ws.Range("D4").Value = "Difference in gas consumption."
$ws.Range("D1").Value = "Value Description"
$ws.Range("D5").Value = "Difference in the amount of produced waste."
$ws.Range("D6").Value = "Euclidean distance difference between the sensors."
$ws.Range("D7").Value = "Difference in the number of parks located around the sensors."
$ws.Range("D8").Value = "Inverse of the total distance between a sensor and all waste facilities."
$ws.Range("A9").Value = "Pollution"
$ws.Range("D9").Value = "Difference in common pollution production."

# 3. Fill in the new PM10 MSE/RMSE result values for the existing adjacency
#    matrix rows (2-7) and the two brand-new adjacency matrix rows (8-9).
$ws.Range("B2").Value = 635.27724244512399
$ws.Range("C2").Value = 25.2047067518176

$ws.Range("B3").Value = 638.65099999999995
$ws.Range("C3").Value = 25.271999999999998

$ws.Range("B4").Value = 657.66272691585402
$ws.Range("C4").Value = 25.644935697245401

$ws.Range("B5").Value = 633.029
$ws.Range("C5").Value = 25.16

$ws.Range("B6").Value = 647.49599999999998
$ws.Range("C6").Value = 25.446000000000002

$ws.Range("B7").Value = 607.89800000000002
$ws.Range("C7").Value = 24.655999999999999

$ws.Range("B8").Value = 615.83057262393595
$ws.Range("C8").Value = 24.815933845494001

$ws.Range("B9").Value = 638.18499999999995
$ws.Range("C9").Value = 25.262

# 4. Number columns (B & C), across the header and every data row, are
#    center-aligned both horizontally and vertically.
$numRange = $ws.Range("B1:C9")
$numRange.HorizontalAlignment = -4108
$numRange.VerticalAlignment = -4108

$ws.Range("A1:D1").HorizontalAlignment = -4108
$ws.Range("A1:D1").VerticalAlignment = -4108

# 5. The two freshly-added rows (8 & 9) need the same look as the existing
#    rows: column A bold/14pt, column D italic/12pt, both left+vcenter.
$newLabelRange = $ws.Range("A8:A9")
$newLabelRange.Font.Bold = $true
$newLabelRange.Font.Size = 14
$newLabelRange.HorizontalAlignment = -4131
$newLabelRange.VerticalAlignment = -4108

$newDescRange = $ws.Range("D8:D9")
$newDescRange.Font.Italic = $true
$newDescRange.Font.Size = 12
$newDescRange.HorizontalAlignment = -4131
$newDescRange.VerticalAlignment = -4108

# 6. Extend the Excel Table ("Tablo13") so it covers the two new rows.
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:D9"))
